$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.517.75'
$ws.Range('E2').Value = '  -2.16%  '
$ws.Range('D3').Value = '2.488.34'
$ws.Range('E3').Value = '  -1.14%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.39'
$ws.Range('E5').Value = '  +0.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.94'
$ws.Range('E6').Value = '  -3.70%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.550'
$ws.Range('E7').Value = '  -2.36%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.501'
$ws.Range('E9').Value = '  -2.90%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '33.67'
$ws.Range('E10').Value = '  -4.22%  '
$ws.Range('E11').Value = '  -2.41%  '
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.02'
$ws.Range('E13').Value = '  -2.51%  '
$ws.Range('D14').Value = '2.872.18'
$ws.Range('E14').Value = '  -1.10%  '
$ws.Range('E15').Value = '  +1.39%  '
$ws.Range('D16').Value = '2.481.61'
$ws.Range('E16').Value = '  +0.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.795'
$ws.Range('E17').Value = '  -1.19%  '
$ws.Range('D18').Value = '41.468.68'
$ws.Range('E18').Value = '  -2.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.35'
$ws.Range('E19').Value = '  -3.64%  '
$ws.Range('E20').Value = '  -1.96%  '
$ws.Range('E21').Value = '  -6.63%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.06'
$ws.Range('E22').Value = '  -0.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.57'
$ws.Range('E23').Value = '  -1.07%  '
$ws.Range('E24').Value = '  -2.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.92'
$ws.Range('E25').Value = '  -3.89%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.30'
$ws.Range('E27').Value = '  -4.06%  '
$ws.Range('E28').Value = '  -0.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.78'
$ws.Range('E29').Value = '  -2.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.54'
$ws.Range('E30').Value = '  -4.16%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '152.68'
$ws.Range('E31').Value = '  -2.64%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.51'
$ws.Range('E32').Value = '  -4.92%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.58'
$ws.Range('E33').Value = '  -3.53%  '
$ws.Range('E34').Value = '  -5.78%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0756'
$ws.Range('E35').Value = '  -3.55%  '
$ws.Range('B36').Value = 'Celestia'
$ws.Range('C36').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.02'
$ws.Range('E36').Value = '  +2.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.09'
$ws.Range('E37').Value = '  -1.59%  '
$ws.Range('E38').Value = '  -3.08%  '
$ws.Range('E39').Value = '  -1.56%  '
$ws.Range('E40').Value = '  -7.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.20'
$ws.Range('E41').Value = '  +1.95%  '
$ws.Range('E42').Value = '  +0.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '19.77'
$ws.Range('E43').Value = '  -9.16%  '
$ws.Range('D44').Value = '2.003.39'
$ws.Range('E44').Value = '  +0.81%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0288'
$ws.Range('E45').Value = '  -2.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.05'
$ws.Range('E46').Value = '  -6.54%  '
$ws.Range('E47').Value = '  -1.74%  '
$ws.Range('D48').Value = '2.733.10'
$ws.Range('E48').Value = '  -1.11%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '70.13'
$ws.Range('E49').Value = '  -1.44%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '97.35'
$ws.Range('E50').Value = '  -3.00%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.178'
$ws.Range('E51').Value = '  -5.55%  '
